$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-18 Tuesday" "2025-11-19 Wednesday"

Replace-Text "70×75=5250" "14×60=840"
Replace-Text "41×59=2419" "51×45=2295"
Replace-Text "61×20=1220" "65×45=2925"
Replace-Text "69×90=6210" "64×60=3840"
Replace-Text "71×38=2698" "77×19=1463"

Replace-Text "95×58=5510" "45×15=675"
Replace-Text "11×83=913" "52×47=2444"
Replace-Text "88×73=6424" "96×51=4896"
Replace-Text "46×36=1656" "50×76=3800"
Replace-Text "93×19=1767" "30×81=2430"

Replace-Text "78×21=1638" "48×68=3264"
Replace-Text "23×69=1587" "50×68=3400"
Replace-Text "51×60=3060" "22×43=946"
Replace-Text "59×40=2360" "84×22=1848"
Replace-Text "88×97=8536" "14×50=700"

Replace-Text "62×92=5704" "17×63=1071"
Replace-Text "85×71=6035" "94×80=7520"
Replace-Text "26×63=1638" "53×84=4452"
Replace-Text "51×71=3621" "90×86=7740"
Replace-Text "28×25=700" "98×38=3724"

Replace-Text "20×54=1080" "22×59=1298"
Replace-Text "91×90=8190" "44×93=4092"
Replace-Text "16×53=848" "29×34=986"
Replace-Text "22×89=1958" "98×32=3136"
Replace-Text "93×88=8184" "46×87=4002"
